# Add a new "TA" (Tamil) sheet between "HI" and "EN" holding the
# "Gnana Nirmalame" bhajan entry, and tidy up the selection left on "HI".

$wb = $excel.ActiveWorkbook

$hi = $wb.Worksheets.Item("HI")
$en = $wb.Worksheets.Item("EN")

# The HI sheet was left with a stray "C3" selection; set it back to the
# header row selection like the other language sheets.
$hi.Range("A1:D1").Select() | Out-Null

# Insert the new "TA" worksheet right after "HI" (i.e. right before "EN").
$ta = $wb.Worksheets.Add($null, $hi)
$ta.Name = "TA"

# Header row (same column headings as the other language sheets).
$ta.Range("A1").Value = "Sr. No."
$ta.Range("B1").Value = "Code"
$ta.Range("C1").Value = "Name (Roman)"
$ta.Range("D1").Value = "Name (Orig)"
$ta.Range("A1:D1").Font.Bold = $true

# Data row: "Gnana Nirmalame".
$ta.Range("A2").Value = 1
$ta.Range("D2").Value = "ஞான நிர்மலாமே"
$ta.Range("C2").Value = "Gnana Nirmalame"
$ta.Range("B2").Value = "GNNM"

# Leave the selection where Excel would land after typing the row.
$ta.Range("A3").Select() | Out-Null
